# Add "Hint" (H) and "Popup" (I) columns to the header row, matching the
# look of the existing "Query"/"Output" header cells (F1:G1): white bold-ish
# Calibri text on the dark-blue fill, no border, general number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value = "Hint"
$ws.Range("I1").Value = "Popup"

# Clone the visual format of the existing header cells (F1:G1) onto the new
# ones, then strip the parts that should NOT carry over (border, text
# alignment/number format used for the data column underneath) so the new
# header cells end up with just the shared font + fill look.
$ws.Range("F1:G1").Copy()
$ws.Range("H1:I1").PasteSpecial(-4122) | Out-Null

$newHeader = $ws.Range("H1:I1")
$newHeader.Borders.LineStyle = -4142
$newHeader.NumberFormat = "General"
$newHeader.HorizontalAlignment = 1
$newHeader.VerticalAlignment = -4107
$newHeader.WrapText = $false

# Reflect the extended used range in the active selection, as a user would
# leave it positioned under the newly added "Hint" column.
$ws.Range("H8").Select() | Out-Null
